$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$idx = 0
$newValues = @(
  "17+25=",
  "18+54=",
  "59+17=",
  "90-89=",
  "51-46=",
  "91-69=",
  "90-12=",
  "27-8=",
  "5+7=",
  "61-26=",
  "81-74=",
  "41-7=",
  "67-29=",
  "23+69=",
  "81-16=",
  "7+79=",
  "27+57=",
  "75-39=",
  "23-7=",
  "5+16=",
  "43+39=",
  "9+84=",
  "64-25=",
  "91-8=",
  "71-9=",
  "82-66=",
  "29+68=",
  "87-18=",
  "4+47=",
  "62-5=",
  "94-57=",
  "48+28=",
  "53+28=",
  "59+22=",
  "71-29=",
  "91-55=",
  "83-15=",
  "68+9=",
  "80-54=",
  "18+36=",
  "91-25=",
  "70-52=",
  "74-36=",
  "9+42=",
  "17+77=",
  "40-19=",
  "71-44=",
  "30-6=",
  "27+17=",
  "18+44=",
  "24-17=",
  "25+49=",
  "3+88=",
  "15+37=",
  "38+55=",
  "6+79=",
  "8+47=",
  "80-49=",
  "5+9=",
  "71-64=",
  "68-9=",
  "83-37=",
  "38+57=",
  "61-6=",
  "80-71=",
  "16+16=",
  "71-15=",
  "90-24=",
  "50-48=",
  "22+19=",
  "16+39=",
  "18+7=",
  "61-6=",
  "7+87=",
  "92-13=",
  "79+12=",
  "37+38=",
  "73-14=",
  "2+49=",
  "96-77=",
  "51-19=",
  "17+38=",
  "8+27=",
  "91-87=",
  "37-18=",
  "26+65=",
  "23-15=",
  "37+24=",
  "74-9=",
  "43+18=",
  "90-33=",
  "5+48=",
  "33-6=",
  "5+19=",
  "39+59=",
  "6+16=",
  "36-18=",
  "90-53=",
  "7+5=",
  "7+59="
)
for ($row = 1; $row -le 20; $row++) {
  for ($col = 1; $col -le 5; $col++) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newValues[$idx]
    $idx++
  }
}
